$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume/1h (E) values for existing rows ---
# D holds text like '583.71' or '68.041.63' -- force Text format so
# Excel doesn't silently coerce it to a Double (losing formatting).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.041.63'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.336.35'
$ws.Range("E3").Value = '  +0.50%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.71'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.15'
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("E8").Value = '  +1.93%  '
$ws.Range("E9").Value = '  +3.85%  '
$ws.Range("E10").Value = '  +1.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '47.85'
$ws.Range("E11").Value = '  +5.17%  '
$ws.Range("E12").Value = '  +1.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '701.61'
$ws.Range("E13").Value = '  +4.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.875.97'
$ws.Range("E14").Value = '  +0.69%  '
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.083.75'
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.343.14'
$ws.Range("E18").Value = '  +1.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.48'
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.14'
$ws.Range("E20").Value = '  +2.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.895'
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.40'
$ws.Range("E22").Value = '  +0.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '16.97'
$ws.Range("E23").Value = '  -1.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '100.63'
$ws.Range("E24").Value = '  +3.05%  '
$ws.Range("E25").Value = '  +1.59%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.43'
$ws.Range("E27").Value = '  +2.03%  '
$ws.Range("E28").Value = '  -0.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.54'
$ws.Range("E29").Value = '  +1.26%  '
$ws.Range("E30").Value = '  -2.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '572.16'
$ws.Range("E31").Value = '  -2.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.99'
$ws.Range("E32").Value = '  +0.30%  '
$ws.Range("E33").Value = '  +1.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.745.88'
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.13'
$ws.Range("E36").Value = '  +3.12%  '
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.45'
$ws.Range("E38").Value = '  +9.60%  '
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("E42").Value = '  +1.48%  '
$ws.Range("E43").Value = '  +0.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.31'
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("E45").Value = '  +0.37%  '
$ws.Range("E46").Value = '  +0.91%  '
$ws.Range("E47").Value = '  +1.12%  '
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("E49").Value = '  -1.29%  '

# --- Rows reordered with refreshed data (coin swapped position) ---
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.63'
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.15'
$ws.Range("E41").Value = '  +0.91%  '
$ws.Range("B50").Value = 'CoreDAO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.76'
$ws.Range("E50").Value = '  +4.36%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '130.19'
$ws.Range("E51").Value = '  +0.50%  '
